# Apply the BOM update described in the commit:
# "added level shifter and button header. Need antenna"
#
# 1. Remove the (unused) external workbook reference/link.
# 2. Add a new row (14) to the BOM for a level shifter part.
# 3. Move the active selection to E15 (next empty cost cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove external references (control_board_3.1_parts.xlsx link) ---
$links = $wb.LinkSources()
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# --- Add new BOM row for the level shifter part ---
$ws.Range("A14").Value = "Level shifter"
$ws.Range("B14").Value = "CD40109BPWR"
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0.43

# --- Update the active selection to match the edited workbook state ---
[void]$ws.Range("E15").Select()
